$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 347-348, shifting the existing data (old rows 347-398)
# down to rows 349-400. This mirrors the weekly roll of the dataset: a new
# week's worth of data is prepended and the oldest week's rows are preserved
# at the bottom.
$ws.Rows("347:348").Insert()

# New row 347 (Primera)
$ws.Range("A347").Value = 8
$ws.Range("B347").Value = "Terminal La Palmera de La Serena"
$ws.Range("C347").Value = "Coquimbo"
$ws.Range("D347").Value = 44474
$ws.Range("E347").Value = 4
$ws.Range("F347").Value = 100112008
$ws.Range("G347").Value = "Coliflor"
$ws.Range("H347").Value = "Sin especificar"
$ws.Range("I347").Value = "Primera"
$ws.Range("J347").Value = 2200
$ws.Range("K347").Value = 650
$ws.Range("L347").Value = 700
$ws.Range("M347").Value = 675
$ws.Range("N347").Value = '$/unidad'
$ws.Range("O347").Value = "Provincia del Elquí"
$ws.Range("P347").Value = 675
$ws.Range("Q347").Value = 1
$ws.Range("R347").Value = "Hortaliza"

# New row 348 (Segunda)
$ws.Range("A348").Value = 8
$ws.Range("B348").Value = "Terminal La Palmera de La Serena"
$ws.Range("C348").Value = "Coquimbo"
$ws.Range("D348").Value = 44474
$ws.Range("E348").Value = 4
$ws.Range("F348").Value = 100112008
$ws.Range("G348").Value = "Coliflor"
$ws.Range("H348").Value = "Sin especificar"
$ws.Range("I348").Value = "Segunda"
$ws.Range("J348").Value = 1360
$ws.Range("K348").Value = 550
$ws.Range("L348").Value = 600
$ws.Range("M348").Value = 575
$ws.Range("N348").Value = '$/unidad'
$ws.Range("O348").Value = "Provincia del Elquí"
$ws.Range("P348").Value = 575
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = "Hortaliza"
